$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.137.15'
$ws.Range('E2').Value = '  +0.68%  '
$ws.Range('D3').Value = '2.796.90'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '361.96'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '110.31'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.51%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.563'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +3.65%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.596'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.37'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0856'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.79%  '
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.65'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.70'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').Value = '3.225.69'
$ws.Range('E15').Value = '  -1.00%  '
$ws.Range('D16').Value = '2.774.19'
$ws.Range('E16').Value = '  -1.79%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.931'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +5.15%  '
$ws.Range('D18').Value = '51.990.24'
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.43'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.13'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.25'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.53%  '
$ws.Range('D22').Value = '0.0₃0986'
$ws.Range('E22').Value = '  +1.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '273.68'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.87'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.78'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.73'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.68%  '
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.23'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.91%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.23'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.143'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0473'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +8.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '51.28'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.41'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.78'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.45'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +12.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0842'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.998'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.31%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.23'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.38%  '
$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.29'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.46%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.01'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.93%  '
$ws.Range('E41').Value = '  +3.66%  '
$ws.Range('E42').Value = '  +0.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '125.42'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.24'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.42%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.96'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.18%  '
$ws.Range('D46').Value = '2.068.29'
$ws.Range('E46').Value = '  +1.55%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.26'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.19%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.31'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.75'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.80%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.952'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.05'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.92%  '
